$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 19 (the PVL row), pushing PVL down to row 20
$ws.Rows.Item(19).Insert()

# Match the row height/formatting used by the rest of the data rows
$ws.Rows.Item(19).RowHeight = $ws.Rows.Item(18).RowHeight

# Copy the formatting (style) used by the other data rows onto the new row
$ws.Range("A18:D18").Copy()
$ws.Range("A19:D19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row 19 with the "Unknown" cell type entry (no rationale text)
$ws.Range("A19").Value = "Unknown"
$ws.Range("B19").Value = "#808080"
$ws.Range("C19").Value = "Gray"
$ws.Range("D19").Value = ""
